$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 20 new "Entrainement" rows (J+2 on 2025-10-14 and J-3 on 2025-10-15),
# columns A..V, appended after the existing last row (677).
$newRows = @(
    @('Entrainement',45944,'Global','J+2','Karahali Souaré','right forward','01:18:18',5.64,0.15,5.48,0.13,0.03,0.01,0,1,3.93,26.91,5.49,44,11,36,10),
    @('Entrainement',45944,'Global','J+2','Emmanuel Valey','left forward','01:20:33',6.02,0.09,5.93,0.08,0.01,0,0,0,4.15,24.05,5.13,48,11,17,2),
    @('Entrainement',45944,'Global','J+2','Romain Thunet','center back','01:20:19',5.79,0.02,5.77,0.02,0,0,0,0,4.3,16.27,3.53,10,0,4,0),
    @('Entrainement',45944,'Global','J+2','Omar Benyounes','center midfield','01:17:51',5.65,0.07,5.57,0.07,0,0,0,0,4,19.4,3.95,13,0,6,2),
    @('Entrainement',45944,'Global','J+2','Jeremie Laurent','left forward','01:20:26',6.13,0.21,5.91,0.18,0.04,0,0,0,4.5,23.2,4.84,46,8,19,12),
    @('Entrainement',45944,'Global','J+2','Mattheo Haon','right back','01:21:01',6.76,0.42,6.32,0.35,0.08,0,0,0,4.92,23.42,4.61,41,6,14,5),
    @('Entrainement',45944,'Global','J+2','Malik Boussaid','right back','01:20:05',6.76,0.27,6.48,0.27,0.01,0,0,0,4.56,23.43,4.51,40,4,32,5),
    @('Entrainement',45945,'Global','J-3','Mattheo Haon','right back','01:43:19',8.19,1.18,6.99,0.68,0.4,0.12,0,12,4.7,27.6,4.7,36,9,27,10),
    @('Entrainement',45945,'Global','J-3','Yoann Martelat','center midfield','01:45:02',7.78,0.96,6.81,0.62,0.24,0.1,0,7,4.42,27.66,4.49,12,4,11,3),
    @('Entrainement',45945,'Global','J-3','Amine Taiar','center back','01:42:31',7.2,0.97,6.22,0.47,0.38,0.13,0,9,4.14,28.72,4.17,44,1,22,6),
    @('Entrainement',45945,'Global','J-3','Kamal Bafounta','center midfield','01:43:19',8.07,1.26,6.79,0.91,0.33,0.03,0,6,4.61,26.81,4.28,29,3,20,7),
    @('Entrainement',45945,'Global','J-3','Romain Thunet','center back','01:42:41',7.47,0.97,6.49,0.58,0.31,0.09,0,7,4.36,27.4,4.27,33,3,30,7),
    @('Entrainement',45945,'Global','J-3','Levy Ndoutoume','left back','01:43:09',6.96,1.21,5.73,0.68,0.36,0.18,0,13,3.99,30.05,4.62,46,10,30,8),
    @('Entrainement',45945,'Global','J-3','Ilan Ihaddadene','center midfield','01:42:49',8.65,1.14,7.49,0.75,0.29,0.11,0.01,7,4.97,31.28,5.06,41,15,25,6),
    @('Entrainement',45945,'Global','J-3','Omar Benyounes','center midfield','01:44:34',8.35,1.19,7.15,0.71,0.29,0.21,0,12,4.74,29.25,4.65,38,12,32,10),
    @('Entrainement',45945,'Global','J-3','Malik Boussaid','right back','01:38:36',8.34,1.37,6.94,0.84,0.35,0.16,0.05,14,4.71,32.66,5.19,87,16,61,21),
    @('Entrainement',45945,'Global','J-3','Emmanuel Valey','left forward','00:45:37',4.89,1,3.87,0.66,0.27,0.09,0,5,5.84,30.15,5.11,36,14,28,7),
    @('Entrainement',45945,'Global','J-3','Amir Kherrab','center midfield','01:43:47',7.12,0.82,6.29,0.55,0.22,0.07,0,4,4.02,28.91,4.03,20,1,15,4),
    @('Entrainement',45945,'Global','J-3','Karahali Souaré','right forward','01:43:20',8.05,1.09,6.94,0.61,0.34,0.15,0.01,14,4.36,31.81,5.8,74,23,52,23),
    @('Entrainement',45945,'Global','J-3','Jeremie Laurent','left forward','01:43:46',8.15,1.44,6.69,0.76,0.42,0.27,0,12,4.65,30.11,5.06,54,23,35,12)

)

$startRow = 678
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

$endRow = $startRow + $newRows.Length - 1

# Match the date-format style used throughout column B (e.g. B2).
$ws.Cells.Item(2, 2).Copy() | Out-Null
$ws.Range($ws.Cells.Item($startRow, 2), $ws.Cells.Item($endRow, 2)).PasteSpecial(-4122) | Out-Null

# Match the centered style used throughout column D (e.g. D650).
$ws.Cells.Item(650, 4).Copy() | Out-Null
$ws.Range($ws.Cells.Item($startRow, 4), $ws.Cells.Item($endRow, 4)).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Move the selection to mirror the saved workbook (cursor sitting just past
# the newly appended data).
$ws.Cells.Item(700, 3).Select() | Out-Null
